$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 2), matching the header row's (B1:P1) bold / bordered /
# centered-top style so A2 picks up the same look the rest of the table uses.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "['', '']"
$ws.Range("C2").Value = "MetaDiff"
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 32
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 256
$ws.Range("K2").Value = 0.9873920511210893
$ws.Range("L2").Value = 0.9812361746565094
$ws.Range("M2").Value = 0.9813206132602249
$ws.Range("N2").Value = 0.9890020384235453
$ws.Range("O2").Value = 404.6989099979401
$ws.Range("P2").Value = 202.3413822650909

# Only A2 carries the header-matching style in the target workbook.
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1
